# Apply updated NOAA temperature data (average_county_temperature, column K)
# and the resulting recalculated worst/best ASHP COP values (columns R and S).
# Facilities 1005513 and 1005849 (rows 10-13) are unaffected by this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  K = 15.74228395061728; R = 1.837513876759573; S = 2.005936573945218 },
    @{ Row = 3;  K = 15.74228395061728; R = 1.837513876759573; S = 2.005936573945218 },
    @{ Row = 4;  K = 13.46442495126706; R = 1.798225615362447; S = 1.958604378795604 },
    @{ Row = 5;  K = 13.46442495126706; R = 1.798225615362447; S = 1.958604378795604 },
    @{ Row = 6;  K = 18.89814814814816; R = 1.894871325212932; S = 2.075424331741031 },
    @{ Row = 7;  K = 18.89814814814816; R = 1.894871325212932; S = 2.075424331741031 },
    @{ Row = 8;  K = 13.46442495126706; R = 1.798225615362447; S = 1.958604378795604 },
    @{ Row = 9;  K = 13.46442495126706; R = 1.798225615362447; S = 1.958604378795604 },
    @{ Row = 14; K = 15.74228395061728; R = 1.837513876759573; S = 2.005936573945218 },
    @{ Row = 15; K = 15.74228395061728; R = 1.837513876759573; S = 2.005936573945218 },
    @{ Row = 16; K = 21.28240740740739; R = 1.940636870984383; S = 2.131200751448103 },
    @{ Row = 17; K = 21.28240740740739; R = 1.940636870984383; S = 2.131200751448103 },
    @{ Row = 18; K = 15.74228395061728; R = 1.837513876759573; S = 2.005936573945218 },
    @{ Row = 19; K = 15.74228395061728; R = 1.837513876759573; S = 2.005936573945218 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 11).Value = $u.K   # column K = average_county_temperature
    $ws.Cells.Item($u.Row, 18).Value = $u.R   # column R = worst_ashp_cop
    $ws.Cells.Item($u.Row, 19).Value = $u.S   # column S = best_ashp_cop
}
